$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet










$ws.Range("B10").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C10").Value = "9146830 - Danúbia Caporusso Bargos"



$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Representation and presentation forms of environmental characteristics; Environmental characterization and its application in a watershed; Environment susceptibilities and vocations determination and environmental susceptibility concept."
$ws.Range("C14").Value = "Representation and presentation forms of environmental characteristics; Environmental characterization and its application in a watershed; Environment susceptibilities and vocations determination and environmental susceptibility concept."

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"
$ws.Rows(15).RowHeight = 120

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Watershed as a unit of study and management of water resources; Environmental characterization of watershed; Morphometric characterization of watershed; Brazilian hydrographic Regions; Brazilian institutional structure and legal frameworks in water resources."
$ws.Range("C16").Value = "Watershed as a unit of study and management of water resources; Environmental characterization of watershed; Morphometric characterization of watershed; Brazilian hydrographic Regions; Brazilian institutional structure and legal frameworks in water resources."

$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows(17).AutoFit()

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C18").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Rows(18).RowHeight = 60

$ws.Range("A19").Value = "Critério:"

$ws.Range("A20").Value = "Norma de recuperação:"

$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows(22).AutoFit()

$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOB1255 -  Hidrologia Aplicada  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1255 -  Hidrologia Aplicada  (Requisito fraco)`n"
$ws.Rows(23).RowHeight = 30

$ws.Rows(24).Delete()
